# Fruta / hortaliza, semanal
# Rearranges the per-row "Fecha" (D) and price block (J,K,L,M,P) values
# across rows 2-23 according to the new weekly ordering, while leaving the
# rest of each row (Mercado, Region, Variedad, etc. - constant across rows
# in this sheet) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (values currently at the source row move to
# the destination row)
$mapping = @{
    2  = 18
    3  = 8
    4  = 17
    5  = 19
    6  = 4
    7  = 16
    8  = 6
    9  = 10
    10 = 23
    11 = 21
    12 = 7
    13 = 11
    14 = 9
    15 = 3
    16 = 20
    17 = 14
    18 = 22
    19 = 2
    20 = 13
    21 = 15
    22 = 12
    23 = 5
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot every source value first (the mapping is a permutation, not a set
# of independent swaps, so all reads must happen before any writes).
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$srcRow").Value2
    }
    $snapshot[$row] = $rowValues
}

# Now write the snapshotted values into their destination rows.
foreach ($row in $snapshot.Keys) {
    $rowValues = $snapshot[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}
